# Fix typo "Posisi" -> "Jabatan" in the committee table header, and
# resize the first two columns ("No" / "N a m a") of that table
# (504 -> 596 twips, 3069 -> 2977 twips; 20 twips = 1 point).

$d = $word.ActiveDocument

# The committee table is the 3-column table whose grid is 504/3069/1984
# twips wide (25.2/153.45/99.2 points) - find it rather than assuming a
# fixed index, then resize its first two columns.
foreach ($t in $d.Tables) {
    if ($t.Columns.Count -eq 3) {
        $col1 = $t.Columns.Item(1)
        $col2 = $t.Columns.Item(2)
        if ([Math]::Abs($col1.Width - 25.2) -lt 0.2 -and [Math]::Abs($col2.Width - 153.45) -lt 0.2) {
            $col1.Width = 29.8
            $col2.Width = 148.85
        }
    }
}

# Correct the misspelled column header.
$d.Content.Find.Execute("Posisi", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jabatan", 2)
